$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.225.34"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "1.605.38"
$ws.Range("E3").Value = "  +0.36%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("E5").Value = "  -0.02%  "
$ws.Range("D6").Value = "'304.46"
$ws.Range("E6").Value = "  +0.87%  "
$ws.Range("D7").Value = "'0.3762"
$ws.Range("E7").Value = "  -0.58%  "
$ws.Range("D8").Value = "'52.54"
$ws.Range("E8").Value = "  +5.44%  "
$ws.Range("D9").Value = "'0.3629"
$ws.Range("E9").Value = "  -0.41%  "
$ws.Range("D10").Value = "'1.274"
$ws.Range("E10").Value = "  +1.30%  "
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").Value = "'0.08154"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("B12").Value = "BinanceUSD"
$ws.Range("C12").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D12").Value = "'1.001"
$ws.Range("E12").Value = "  +0.04%  "
$ws.Range("D13").Value = "'22.91"
$ws.Range("E13").Value = "  +1.79%  "
$ws.Range("D14").Value = "'6.598"
$ws.Range("E14").Value = "  +0.09%  "
$ws.Range("D15").Value = "'7.396"
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").Value = "'0.00001252"
$ws.Range("E16").Value = "  +0.42%  "
$ws.Range("D17").Value = "1.605.31"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("D18").Value = "'93.92"
$ws.Range("E18").Value = "  +2.10%  "
$ws.Range("D19").Value = "'0.06916"
$ws.Range("E19").Value = "  +1.40%  "
$ws.Range("D20").Value = "'18.15"
$ws.Range("E20").Value = "  -0.41%  "
$ws.Range("D21").Value = "'6.539"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "'1.003"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "'12.93"
$ws.Range("E23").Value = "  -1.01%  "
$ws.Range("D24").Value = "23.223.47"
$ws.Range("E24").Value = "  +0.80%  "
$ws.Range("D25").Value = "'2.450"
$ws.Range("E25").Value = "  +3.54%  "
$ws.Range("D26").Value = "'3.063"
$ws.Range("E26").Value = "  +9.06%  "
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").Value = "'150.04"
$ws.Range("E28").Value = "  -0.17%  "
$ws.Range("D29").Value = "'5.286"
$ws.Range("E29").Value = "  +0.91%  "
$ws.Range("D30").Value = "'135.22"
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("D31").Value = "'2.392"
$ws.Range("E31").Value = "  +2.01%  "
$ws.Range("D32").Value = "'6.759"
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("D33").Value = "1.781.44"
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").Value = "'0.9655"
$ws.Range("E34").Value = "  +0.28%  "
$ws.Range("D35").Value = "'0.07505"
$ws.Range("E35").Value = "  -0.92%  "
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "'10.39"
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.02761"
$ws.Range("E37").Value = "  +1.99%  "
$ws.Range("E38").Value = "  -0.14%  "
$ws.Range("D39").Value = "'6.121"
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("D40").Value = "'0.08804"
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("D41").Value = "'1.426"
$ws.Range("E41").Value = "  +4.71%  "
$ws.Range("D42").Value = "'0.7101"
$ws.Range("E42").Value = "  +1.06%  "
$ws.Range("D44").Value = "'15.83"
$ws.Range("E44").Value = "  +3.91%  "
$ws.Range("D45").Value = "'0.6540"
$ws.Range("E45").Value = "  -1.18%  "
$ws.Range("D46").Value = "'2.333"
$ws.Range("E46").Value = "  +2.05%  "
$ws.Range("D47").Value = "'4.008"
$ws.Range("D48").Value = "'133.93"
$ws.Range("E48").Value = "  +1.08%  "
$ws.Range("D49").Value = "'0.07940"
$ws.Range("E49").Value = "  +0.57%  "
$ws.Range("D50").Value = "'1.210"
$ws.Range("E50").Value = "  -0.87%  "
$ws.Range("D51").Value = "'1.192"
